$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 corresponds to file_name = "metrics_sim_with_priors.json"
$ws.Range("C3").Value = 0.3363095238095238
$ws.Range("D3").Value = 0.5952380952380952
$ws.Range("E3").Value = 0.8809523809523809
$ws.Range("G3").Value = 0.9970238095238095
$ws.Range("H3").Value = 0.3251195292386907
$ws.Range("I3").Value = 0.1912755030673621
$ws.Range("J3").Value = 0.2380952380952381
$ws.Range("K3").Value = 624.3095238095239

$ws.Range("Q3").Value = 29
$ws.Range("R3").Value = 93
$ws.Range("S3").Value = 329
$ws.Range("T3").Value = 671
$ws.Range("U3").Value = 1088
$ws.Range("V3").Value = 2354
$ws.Range("W3").Value = 2290
$ws.Range("X3").Value = 2054
$ws.Range("Y3").Value = 1712
$ws.Range("Z3").Value = 1295

$ws.Range("AF3").Value = 0.98783
$ws.Range("AG3").Value = 0.960974
$ws.Range("AH3").Value = 0.861939
$ws.Range("AI3").Value = 0.718422
$ws.Range("AJ3").Value = 0.5434330000000001
